$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, column letter, new value.
# Derived from the authoritative row-by-row diff of the workbook:
# dataset rows 2-35 were reshuffled/updated (Volumen, Precio minimo/maximo/promedio,
# Unidad de comercializacion, Precio $/Kg, Kg o Unidades, and Fecha all move together
# per record); row 23 is untouched.
$edits = @(
    @(2, "D", 44601),
    @(2, "J", 270),
    @(2, "K", 2200),
    @(2, "L", 2500),
    @(2, "M", 2350),
    @(2, "P", 1175),
    @(3, "D", 44789),
    @(3, "K", 1400),
    @(3, "L", 1500),
    @(3, "M", 1450),
    @(3, "P", 725),
    @(4, "D", 44525),
    @(4, "J", 300),
    @(4, "K", 1400),
    @(4, "L", 1500),
    @(4, "M", 1450),
    @(4, "P", 725),
    @(5, "D", 44427),
    @(5, "J", 250),
    @(5, "K", 1300),
    @(5, "L", 1500),
    @(5, "M", 1400),
    @(5, "P", 700),
    @(6, "D", 44253),
    @(6, "K", 1800),
    @(6, "L", 2000),
    @(6, "M", 1900),
    @(6, "P", 950),
    @(7, "D", 44229),
    @(7, "K", 1800),
    @(7, "L", 2000),
    @(7, "M", 1900),
    @(7, "P", 950),
    @(8, "D", 44726),
    @(8, "J", 250),
    @(8, "K", 2500),
    @(8, "L", 2800),
    @(8, "M", 2650),
    @(8, "P", 1325),
    @(9, "D", 44616),
    @(9, "J", 270),
    @(9, "K", 1300),
    @(9, "L", 1500),
    @(9, "M", 1400),
    @(9, "P", 700),
    @(10, "D", 44817),
    @(10, "J", 300),
    @(10, "K", 900),
    @(10, "L", 1000),
    @(10, "M", 950),
    @(10, "P", 475),
    @(11, "D", 44540),
    @(11, "J", 300),
    @(11, "K", 900),
    @(11, "L", 1000),
    @(11, "M", 950),
    @(11, "P", 475),
    @(12, "D", 44365),
    @(12, "J", 200),
    @(12, "K", 1800),
    @(12, "L", 2000),
    @(12, "M", 1900),
    @(12, "P", 950),
    @(13, "D", 44266),
    @(13, "J", 300),
    @(13, "K", 1700),
    @(13, "L", 1800),
    @(13, "M", 1750),
    @(13, "P", 875),
    @(14, "D", 44392),
    @(14, "J", 250),
    @(14, "K", 1800),
    @(14, "L", 2000),
    @(14, "M", 1900),
    @(14, "P", 950),
    @(15, "D", 44435),
    @(15, "J", 300),
    @(15, "K", 900),
    @(15, "L", 1000),
    @(15, "M", 950),
    @(15, "P", 475),
    @(16, "D", 44544),
    @(16, "K", 900),
    @(16, "L", 1000),
    @(16, "M", 950),
    @(16, "P", 475),
    @(17, "D", 44795),
    @(18, "D", 44936),
    @(18, "J", 350),
    @(18, "K", 3000),
    @(18, "L", 3500),
    @(18, "M", 3357),
    @(18, "P", 1678),
    @(19, "D", 44302),
    @(19, "J", 300),
    @(19, "K", 900),
    @(19, "L", 1000),
    @(19, "M", 950),
    @(19, "P", 475),
    @(20, "D", 44390),
    @(20, "K", 2400),
    @(20, "L", 2500),
    @(20, "M", 2450),
    @(20, "P", 1225),
    @(21, "D", 44161),
    @(21, "J", 270),
    @(22, "D", 44243),
    @(22, "J", 250),
    @(22, "K", 1200),
    @(22, "L", 1300),
    @(22, "M", 1250),
    @(22, "P", 625),
    @(24, "D", 44572),
    @(24, "J", 300),
    @(24, "K", 1400),
    @(24, "L", 1500),
    @(24, "M", 1450),
    @(24, "N", "`$/atado 1,5 a 2 kilos"),
    @(24, "P", 725),
    @(24, "Q", 2),
    @(25, "D", 44385),
    @(25, "K", 2400),
    @(25, "L", 2500),
    @(25, "M", 2450),
    @(25, "P", 1225),
    @(26, "D", 44403),
    @(26, "J", 250),
    @(26, "K", 1800),
    @(26, "L", 2000),
    @(26, "M", 1900),
    @(26, "P", 950),
    @(27, "D", 44172),
    @(27, "J", 200),
    @(27, "K", 1300),
    @(27, "L", 1500),
    @(27, "M", 1400),
    @(27, "P", 700),
    @(28, "D", 44257),
    @(28, "J", 500),
    @(28, "K", 1400),
    @(28, "L", 1500),
    @(28, "M", 1450),
    @(28, "P", 725),
    @(29, "D", 44917),
    @(29, "K", 2700),
    @(29, "L", 3000),
    @(29, "M", 2850),
    @(29, "P", 1425),
    @(30, "D", 44363),
    @(30, "K", 2500),
    @(30, "L", 2800),
    @(30, "M", 2650),
    @(30, "P", 1325),
    @(31, "D", 44438),
    @(31, "K", 950),
    @(31, "M", 975),
    @(31, "P", 488),
    @(32, "D", 44181),
    @(32, "J", 200),
    @(32, "K", 1000),
    @(32, "L", 1200),
    @(32, "M", 1100),
    @(32, "N", "`$/atado"),
    @(32, "P", 1100),
    @(32, "Q", 1),
    @(33, "D", 44468),
    @(33, "K", 900),
    @(33, "L", 1000),
    @(33, "M", 950),
    @(33, "P", 475),
    @(34, "D", 44291),
    @(34, "J", 250),
    @(34, "K", 1800),
    @(34, "L", 2000),
    @(34, "M", 1900),
    @(34, "P", 950),
    @(35, "D", 44447),
    @(35, "J", 300),
    @(35, "K", 900),
    @(35, "L", 1000),
    @(35, "M", 950),
    @(35, "P", 475),
)

foreach ($edit in $edits) {
    $targetRow = $edit[0]
    $col = $edit[1]
    $val = $edit[2]
    $ws.Range("$col$targetRow").Value = $val
}
